# Add retrieve cell value from another sheet
#
# On the "Hyperlink" sheet, add a second row that shows a label ("Sum")
# together with the value retrieved from the "Formula" sheet, where cell
# C4 holds the result of =SUM(C2:C3).

$wb = $excel.ActiveWorkbook

$hyperlinkSheet = $wb.Worksheets.Item("Hyperlink")
$formulaSheet   = $wb.Worksheets.Item("Formula")

# Grab the computed SUM(C2:C3) value from the Formula sheet.
$sumValue = $formulaSheet.Range("C4").Value2

# Write the label and the retrieved value into the Hyperlink sheet.
$hyperlinkSheet.Range("A2").Value = "Sum"
$hyperlinkSheet.Range("B2").Value = $sumValue
